$wb = $excel.ActiveWorkbook

# Sheet "Test Cases": D3 Results -> PASS
$ws1 = $wb.Worksheets.Item("Test Cases")
$ws1.Range("D3").Value = "PASS"

# Sheet "ListingOpen": E3:E19 Results -> PASS
$ws2 = $wb.Worksheets.Item("ListingOpen")
for ($r = 3; $r -le 19; $r++) {
    $ws2.Cells.Item($r, 5).Value = "PASS"
}
